$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Text
    if ($val -ne $null -and $val -like "*_f1_*") {
        $cell.Value = $val -replace "_f1_", "_"
    }
}
